$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 15 and 17 currently have "Time Spent" = "1 hour"; change it to "30 minutes"
$ws.Range("C15").Value = "30 minutes"
$ws.Range("C17").Value = "30 minutes"

# Add a new row 18: Develop the logic structure, 30 minutes, Yes, Kyle (on 6/5/2019)
# Copy formatting from row 17 down to row 18 first so the new row matches
# the look (fonts/number format) of the existing data rows.
$ws.Range("A17:E17").Copy() | Out-Null
$ws.Range("A18:E18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A18").Value = (Get-Date -Year 2019 -Month 6 -Day 5 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("B18").Value = "Develop the logic structure"
$ws.Range("C18").Value = "30 minutes"
$ws.Range("D18").Value = "Yes"
$ws.Range("E18").Value = "Kyle"

# Match the existing data rows' row height (15.75pt, explicit/custom height)
$ws.Rows.Item(18).RowHeight = 15.75

# Update the active selection to match the target state
$ws.Range("D22").Select()
